$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to
# text (NumberFormat "@") before assignment, otherwise Excel COM auto-converts
# them to numeric cells instead of preserving them as literal text.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.994.26'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.828.92'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '244.23'
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").Value = '0.6309'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("D7").Value = '0.9979'
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '0.2936'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("D10").Value = '22.87'
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("D11").Value = '0.07707'
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").Value = '1.831.21'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '4.987'
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = '0.6697'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '82.94'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = '0.000009618'
$ws.Range("E16").Value = '  +5.81%  '
$ws.Range("D17").Value = '6.070'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '29.009.80'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '12.55'
$ws.Range("E19").Value = '  +1.85%  '
$ws.Range("D20").Value = '226.73'
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").Value = '0.9977'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").Value = '7.154'
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").Value = '0.9980'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '159.91'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = '0.1404'
$ws.Range("E25").Value = '  +3.80%  '
$ws.Range("D26").Value = '8.528'
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").Value = '17.89'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").Value = '1.495'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '4.117'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").Value = '4.063'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '1.196'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").Value = '0.05369'
$ws.Range("E32").Value = '  +3.57%  '
$ws.Range("D33").Value = '1.857'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").Value = '0.7428'
$ws.Range("E34").Value = '  +1.62%  '
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").Value = '2.653'
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '1.242.42'
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("D38").Value = '2.754'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '0.01786'
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("D40").Value = '6.627'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").Value = '0.9007'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").Value = '0.9987'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '101.51'
$ws.Range("D44").Value = '1.980.34'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '0.00000000124'
$ws.Range("E45").Value = '  +3.79%  '
$ws.Range("D46").Value = '64.78'
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("D47").Value = '0.5103'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").Value = '0.4068'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("D49").Value = '8.992'
$ws.Range("D50").Value = '0.05766'
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.644'
$ws.Range("E51").Value = '  -0.20%  '
